$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jaana")

# Restore the workbook window's on-screen position/size (MainWindow fix)
try {
    $win = $excel.ActiveWindow
    $win.Left = 33450
    $win.Top = 1860
    $win.Width = 21630
    $win.Height = 11310
} catch {
}

# Fill in row 15 with a new time-tracking entry, matching the style of the
# preceding date rows (A7:A14) by copying A14's format down onto A15
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A15").Value = 44986
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = "Scrum daily, retro, review ja planning"

# Move the active cell selection to A16
$ws.Range("A16").Select()
